$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "##name##"
$ws.Range("A2").Value = "John"

$ws.Range("B1").Clear()
$ws.Range("B2").Clear()
